# Apply the PN13 ConceptMap metadata refresh + removal of the
# "Id_prescripteur@Phast-uri_nomenclature" mapping row.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL: hl7.fr/fhir/fr/medication -> hl7.fr/ig/fhir/medication
$meta.Range("B2").Value2 = "https://hl7.fr/ig/fhir/medication/ConceptMap/PN13-FHIR-prescmed-practitioner-id-seul-conceptmap"

# Date refresh
$meta.Range("B8").Value2 = "2026-01-15T08:54:26+00:00"

# Jurisdiction was blank, now set to FRANCE
$meta.Range("B11").Value2 = "FRANCE"

# Target: hl7.fr/fhir/fr/medication -> hl7.fr/ig/fhir/medication
$meta.Range("B16").Value2 = "https://hl7.fr/ig/fhir/medication/StructureDefinition/fr-inpatient-medicationrequest"

# --- Mapping Table 0 sheet ---------------------------------------------
$map = $wb.Worksheets.Item("Mapping Table 0")

# Remove the 4th data row (Id_prescripteur@Phast-uri_nomenclature ->
# MedicationRequest.encounter.identifier.system), shifting nothing below
# it since it was the last row.
$map.Rows.Item(4).Delete()
